$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$newUrl = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/local-race-cd"
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = $newUrl
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet updates ---
$elements = $wb.Worksheets.Item("Elements")

# The Extension.url element's Fixed Value equals the StructureDefinition's own
# canonical URL, so it must be renamed in lockstep with the Metadata URL above.
$elements.Range("Q5").Value = $newUrl

# The "ele-1/ext-1" constraint text moves from the base "Extension" row (row 2)
# down to the "Extension.extension" row (row 4), Constraint(s) column (AI).
$constraintText = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}`next-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"

$elements.Range("AI2").Value = ""
$elements.Range("AI4").Value = $constraintText
